$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A291").Value = "IMX-USD"
$ws.Range("A292").Value = "TAO-USD"
$ws.Range("A293").Value = "MNT-USD"
